# Update F-column (price/view count) values on 展览, 演出, and 全部类型 sheets
# to reflect the refreshed data snapshot, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3557
$ws1.Range("F5").Value = 3557
$ws1.Range("F6").Value = 258
$ws1.Range("F7").Value = 5073
$ws1.Range("F8").Value = 5073
$ws1.Range("F9").Value = 518
$ws1.Range("F10").Value = 352
$ws1.Range("F14").Value = 84
$ws1.Range("F15").Value = 31
$ws1.Range("F16").Value = 695
$ws1.Range("F17").Value = 314
$ws1.Range("F22").Value = 362
$ws1.Range("F23").Value = 4905
$ws1.Range("F24").Value = 4905
$ws1.Range("F28").Value = 6027
$ws1.Range("F29").Value = 20
$ws1.Range("F30").Value = 16
$ws1.Range("F32").Value = 332
$ws1.Range("F33").Value = 710
$ws1.Range("F36").Value = 121
$ws1.Range("F38").Value = 1002
$ws1.Range("F42").Value = 866
$ws1.Range("F43").Value = 974
$ws1.Range("F44").Value = 2024

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 54

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3557
$ws4.Range("F8").Value = 3557
$ws4.Range("F9").Value = 258
$ws4.Range("F10").Value = 5073
$ws4.Range("F11").Value = 5073
$ws4.Range("F12").Value = 518
$ws4.Range("F13").Value = 352
$ws4.Range("F17").Value = 84
$ws4.Range("F18").Value = 31
$ws4.Range("F19").Value = 695
$ws4.Range("F20").Value = 314
$ws4.Range("F26").Value = 362
$ws4.Range("F27").Value = 4905
$ws4.Range("F28").Value = 4905
$ws4.Range("F32").Value = 6027
$ws4.Range("F33").Value = 20
$ws4.Range("F34").Value = 16
$ws4.Range("F36").Value = 332
$ws4.Range("F37").Value = 710
$ws4.Range("F41").Value = 121
$ws4.Range("F43").Value = 1002
$ws4.Range("F47").Value = 866
$ws4.Range("F48").Value = 974
$ws4.Range("F50").Value = 2024
$ws4.Range("F52").Value = 54
